$p = $ppt.ActivePresentation

function Set-DateText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "9/9/2020") {
                $sh.TextFrame.TextRange.Text = "9/11/2020"
            }
        }
    }
}

# Slide master date placeholder
Set-DateText $p.SlideMaster.Shapes

# Slide layout date placeholders
for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    $lay = $p.SlideMaster.CustomLayouts.Item($j)
    Set-DateText $lay.Shapes
}
